$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OP_SZIE")
$ws.Activate()

# Insert a new column before column K (shifts old K->L, old L->M)
$ws.Columns("K:K").Insert()

# New header for the inserted column + literal trigger value used for every row
$ws.Range("K2").Value = "TRIGGER_GC"
$ws.Range("K3:K7").Value = 9

# Match the column width Excel computed for the new "TRIGGER_GC" header (best effort,
# COM ColumnWidth is quantized to whole pixels so this is the closest achievable value)
$ws.Columns("K:K").ColumnWidth = 12.86

# Match the recorded selection state after the edit
$ws.Range("M16").Select()
